# Update the "department" column (C) from "FACULTY OF HOSPITALITY" to more
# specific labels, and clear the "promotionValidity" column (R) since the
# 2021 promotion has expired.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 2-5: individual courses -> "Hospitality"
$ws.Range("C2:C5").Value = "Hospitality"

# Rows 6-8: bundled courses -> "Packages"
$ws.Range("C6:C8").Value = "Packages"

# Clear the outdated promotion validity text for all data rows.
$ws.Range("R2:R8").ClearContents()
